$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# --- Update existing rows 2-10 ---

# Row 2
$ws.Range("B2").Value = "5029547-06.2022.8.21.0022"
$ws.Range("C2").Value = "5007960-93.2020.8.21.0022"
$ws.Range("D2").Value = "Originário"
Set-TextValue $ws.Range("I2") "28/09/2022"

# Row 3
$ws.Range("B3").Value = "5002020-12.2020.8.21.0067"
$ws.Range("C3").Value = "9000539-43.2020.8.21.0067"
$ws.Range("D3").Value = "Migrado"
Set-TextValue $ws.Range("I3") "03/11/2020"

# Row 4
$ws.Range("B4").Value = "5002025-34.2020.8.21.0067"
$ws.Range("C4").Value = "9000434-66.2020.8.21.0067"
$ws.Range("D4").Value = "Migrado"
Set-TextValue $ws.Range("I4") "23/09/2020"

# Row 5
$ws.Range("B5").Value = "5000995-32.2018.8.21.0067"
$ws.Range("C5").Value = "9000796-39.2018.8.21.0067"
$ws.Range("D5").Value = "Migrado"
Set-TextValue $ws.Range("I5") "27/11/2018"

# Row 6
$ws.Range("B6").Value = "5002032-26.2020.8.21.0067"
$ws.Range("C6").Value = "9000450-20.2020.8.21.0067"
$ws.Range("D6").Value = "Migrado"
$ws.Range("E6").Value = "Sem dados de processo originário 2"
$ws.Range("F6").Value = "Nulo"
Set-TextValue $ws.Range("I6") "30/09/2020"

# Row 7
$ws.Range("B7").Value = "5008633-12.2022.8.21.4001"
$ws.Range("C7").Value = "5002299-64.2019.8.21.4001"
$ws.Range("D7").Value = "Originário"
Set-TextValue $ws.Range("I7") "14/10/2022"

# Row 8
$ws.Range("B8").Value = "5000229-27.2011.8.21.0001"
$ws.Range("C8").Value = "0420539-74.2011.8.21.0001"
$ws.Range("D8").Value = "Digitalizado"
Set-TextValue $ws.Range("I8") "07/12/2011"

# Row 9
$ws.Range("B9").Value = "5000297-74.2011.8.21.0001"
$ws.Range("C9").Value = "0413928-08.2011.8.21.0001"
Set-TextValue $ws.Range("I9") "02/12/2011"

# Row 10
$ws.Range("B10").Value = "5000320-49.2013.8.21.0001"
$ws.Range("C10").Value = "0044684-94.2013.8.21.0001"
Set-TextValue $ws.Range("I10") "20/02/2013"

# --- Add new rows 11-19 ---

$newRows = @(
    @{ Row=11; A=9;  B="5003815-38.2012.8.21.0001"; C="0023925-46.2012.8.21.0001"; D="Digitalizado"; I="24/01/2012" },
    @{ Row=12; A=10; B="5002031-41.2020.8.21.0067"; C="9000380-03.2020.8.21.0067"; D="Migrado";      I="14/09/2020" },
    @{ Row=13; A=11; B="5002033-11.2020.8.21.0067"; C="9000365-34.2020.8.21.0067"; D="Migrado";      I="09/09/2020" },
    @{ Row=14; A=12; B="5002037-48.2020.8.21.0067"; C="9000624-29.2020.8.21.0067"; D="Migrado";      I="13/11/2020" },
    @{ Row=15; A=13; B="5002090-63.2019.8.21.0067"; C="9001254-22.2019.8.21.0067"; D="Migrado";      I="11/12/2019" },
    @{ Row=16; A=14; B="5002030-56.2020.8.21.0067"; C="9000405-16.2020.8.21.0067"; D="Migrado";      I="18/09/2020" },
    @{ Row=17; A=15; B="5002091-48.2019.8.21.0067"; C="9001270-73.2019.8.21.0067"; D="Migrado";      I="17/12/2019" },
    @{ Row=18; A=16; B="5006127-55.2022.8.21.0059"; C="5001247-88.2020.8.21.0059"; D="Originário";   I="19/09/2022" },
    @{ Row=19; A=17; B="5092615-95.2019.8.21.0001"; C="9015530-62.2019.8.21.0001"; D="Migrado";      I="04/04/2019" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy formatting from the A-column cell of the row above, so style "1"
    # (bold, border, centered) carries over to the newly created row.
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = "Sem dados de processo originário 2"
    $ws.Cells.Item($row, 6).Value = "Nulo"
    $ws.Cells.Item($row, 7).Value = "Sem dados de processo originário 3"
    $ws.Cells.Item($row, 8).Value = "Nulo"
    Set-TextValue $ws.Cells.Item($row, 9) $r.I
}

$excel.CutCopyMode = 0
